$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns: K = Sprint Review, L = Sprint Retrospective ---
$ws.Range("K1").Value = "Sprint Review"
$ws.Range("L1").Value = "Sprint Retrospective"

$ws.Range("L2").Value = "Inga förändringar"
$ws.Range("K2").Value = "Funkar som förväntat"

$ws.Range("L3").Value = "Inga förändringar"
$ws.Range("K3").Value = "Funkar som förväntat"

$ws.Range("L4").Value = "Inga förändringar"
$ws.Range("K4").Value = "Funkar som förväntat"

# --- Header formatting (bold + gray fill), matching the look of the other headers ---
$headerRange = $ws.Range("K1:L1")
$headerRange.Font.Bold = $true
$headerRange.Font.Name = "Arial"
$headerRange.Font.Size = 10
$headerRange.Interior.Pattern = 1
$headerRange.Interior.ThemeColor = 0
$headerRange.Interior.TintAndShade = -0.249977111117893
$headerRange.WrapText = $true



# --- Column width for L (target stored width 22.28515625) ---
$ws.Columns("L").ColumnWidth = 21.451822916666668

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ---
$ws.Range("K4").Select()
